$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G ("discount") values were changed; some cleared, others updated. ---
$ws.Range("G3").Value = 80
$ws.Range("G4").Value = 10
$ws.Range("G5").Value = 10
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("G12").Value = 50
$ws.Range("G13").Value = 100
$ws.Range("G14").Value = 50
$ws.Range("G16").ClearContents()
$ws.Range("G17").Value = 100
$ws.Range("G18").Value = 100

# --- Column G is now formatted as a whole number (was General / 0.00). ---
$ws.Range("G1:G18").NumberFormat = "0"

# --- Selection moved. ---
$ws.Range("H9").Select()
